# Daily scrape update - 2025-07-31 03:45:12 UTC
# Refresh the opportunity listing rows 2-11 with newly scraped data,
# drop the two oldest rows (previously rows 12-13), reset the
# "PREMIUM" highlight on the two rows that no longer qualify, and
# tweak a few column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two trailing rows (old rows 12 & 13) -----------------
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# Opportunity IDs in column A are plain digit strings that must stay
# TEXT (as in the source file) rather than be auto-coerced to numbers.
# A leading apostrophe forces text entry like a user typing in Excel;
# resetting the style back to Normal afterwards drops the transient
# "quote prefix" formatting so the cell ends up unstyled, same as the
# source.
function Set-TextId($cell, $id) {
    $cell.Value = "'" + $id
    $cell.Style = "Normal"
}

# --- Row 2 -------------------------------------------------------------
Set-TextId $ws.Range("A2") "1326636"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326636"
$ws.Range("C2").Value = "Computer and AI Coordinator"
$ws.Range("D2").Value = "London, UK"
$ws.Range("E2").Value = "No"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Capital Care Homes"

# --- Row 3 -------------------------------------------------------------
Set-TextId $ws.Range("A3") "1326533"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326533"
$ws.Range("F3").Value = "4 applicants"

# --- Row 4 -------------------------------------------------------------
Set-TextId $ws.Range("A4") "1326468"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326468"
$ws.Range("C4").Value = "Social Media Executive"
$ws.Range("D4").Value = "Navi Mumbai, Maharashtra, India"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Information Technology Learning Hub (ITLH)"

# --- Row 5 -------------------------------------------------------------
Set-TextId $ws.Range("A5") "1326068"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1326068"
$ws.Range("C5").Value = "Supply Chain Intern (Graduate Internship Program: GRIP)"
$ws.Range("D5").Value = "Nairobi, Kenya"
$ws.Range("F5").Value = "15 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "World Vision"

# --- Row 6 -------------------------------------------------------------
Set-TextId $ws.Range("A6") "1325954"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1325954"
$ws.Range("C6").Value = "Intern; Food Initiative / Africa Regional Office"
$ws.Range("D6").Value = "Nairobi, Kenya"
$ws.Range("F6").Value = "16 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Rockefeller Foundation - Kenya"

# --- Row 7 -------------------------------------------------------------
Set-TextId $ws.Range("A7") "1325634"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1325634"
$ws.Range("C7").Value = "Technical Presales Engineer"
$ws.Range("D7").Value = "Nairobi, Kenya"
$ws.Range("E7").Value = "No"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "9 applicants"
$ws.Range("H7").Value = "HUAWEI"

# --- Row 8 -------------------------------------------------------------
Set-TextId $ws.Range("A8") "1325633"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1325633"
$ws.Range("C8").Value = "IT Infrastructure Operations Engineer"
$ws.Range("D8").Value = "Nairobi, Kenya"
$ws.Range("F8").Value = "10 applicants"
$ws.Range("H8").Value = "HUAWEI"

# --- Row 9 -------------------------------------------------------------
Set-TextId $ws.Range("A9") "1325541"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1325541"
$ws.Range("C9").Value = "Web developer"
$ws.Range("D9").Value = "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"
$ws.Range("F9").Value = "28 applicants"
$ws.Range("H9").Value = "breem_solutions"

# --- Row 10 ------------------------------------------------------------
Set-TextId $ws.Range("A10") "1322500"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1322500"
$ws.Range("C10").Value = "Accelerate Romania | Web Design & Supplier Relations Assistant"
$ws.Range("D10").Value = "Iași, Romania"
$ws.Range("F10").Value = "33 applicants"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "Kombu"

# --- Row 11 ------------------------------------------------------------
Set-TextId $ws.Range("A11") "1316788"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1316788"
$ws.Range("C11").Value = "Travel Coordinator"
$ws.Range("D11").Value = "Mexico City, CDMX, Mexico"
$ws.Range("F11").Value = "79 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "Ikan Experience"

# --- Column width tweaks ------------------------------------------------
# ColumnWidth is offset from the stored OOXML "width" by 5/6 of a
# character (the built-in cell padding), so subtract that to land on
# the exact target width.
$ws.Columns.Item(3).ColumnWidth = 64.16666666666667   # 85 -> 65
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666  # 17 -> 16
$ws.Columns.Item(8).ColumnWidth = 44.166666666666664  # 32 -> 45
